# Update cryptocurrency price (column D) and volume/change (column E) values
# for rows 2-51 to reflect the refreshed data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.573.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.16%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.651.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.10%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.29%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.63%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.48%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3247"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.70%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07021"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.36%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.28%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.987"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.653.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.13%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.587"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.30%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001038"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.98%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06610"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.54%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.23%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.22%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.896"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.55%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.537.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.441"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.61%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.364"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -14.87%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.75%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.839.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.97%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.199"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.46%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.54%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.069"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.773"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -14.80%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08450"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.56%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.673"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.51%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -11.36%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.268"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.42%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.180"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.49%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05996"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.82%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02221"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.61%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2065"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.39%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.133"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -11.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.26%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5894"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.29%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.832"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.29%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5592"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.91%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.45%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.940"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.06%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06933"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.18%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.191"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.56%  "

